$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.782.23"
$ws.Range("E2").Value = "  -0.60%  "
$ws.Range("D3").Value = "'2.299.20"
$ws.Range("E3").Value = "  -0.09%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'306.84"
$ws.Range("E5").Value = "  +2.20%  "
$ws.Range("D6").Value = "'96.16"
$ws.Range("E6").Value = "  -2.13%  "
$ws.Range("E7").Value = "  -2.53%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  -2.75%  "
$ws.Range("D10").Value = "'35.28"
$ws.Range("E10").Value = "  -2.55%  "
$ws.Range("D11").Value = "'0.0791"
$ws.Range("E11").Value = "  +0.16%  "
$ws.Range("D12").Value = "'18.56"
$ws.Range("E12").Value = "  +4.65%  "
$ws.Range("E13").Value = "  +1.30%  "
$ws.Range("E14").Value = "  -1.77%  "
$ws.Range("D15").Value = "'2.658.71"
$ws.Range("E15").Value = "  +0.01%  "
$ws.Range("D16").Value = "'2.304.56"
$ws.Range("E16").Value = "  +1.64%  "
$ws.Range("E17").Value = "  -0.82%  "
$ws.Range("D18").Value = "'42.704.00"
$ws.Range("E18").Value = "  -0.55%  "
$ws.Range("D19").Value = "'13.01"
$ws.Range("E19").Value = "  +1.86%  "
$ws.Range("E20").Value = "  -1.51%  "
$ws.Range("E21").Value = "  -1.10%  "
$ws.Range("D22").Value = "'67.26"
$ws.Range("E22").Value = "  -2.43%  "
$ws.Range("D23").Value = "'235.88"
$ws.Range("E23").Value = "  -0.69%  "
$ws.Range("D24").Value = "'2.13"
$ws.Range("E24").Value = "  -0.57%  "
$ws.Range("E25").Value = "  +0.81%  "
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("E27").Value = "  +0.15%  "
$ws.Range("D28").Value = "'25.14"
$ws.Range("E28").Value = "  +0.91%  "
$ws.Range("E29").Value = "  +16.82%  "
$ws.Range("D30").Value = "'166.58"
$ws.Range("E30").Value = "  +0.89%  "
$ws.Range("E31").Value = "  -0.69%  "
$ws.Range("D32").Value = "'33.02"
$ws.Range("E32").Value = "  +0.02%  "
$ws.Range("E33").Value = "  +0.10%  "
$ws.Range("D34").Value = "'4.75"
$ws.Range("E34").Value = "  -0.44%  "
$ws.Range("D35").Value = "'4.98"
$ws.Range("E35").Value = "  -2.00%  "
$ws.Range("D36").Value = "'17.67"
$ws.Range("E36").Value = "  -1.15%  "
$ws.Range("E37").Value = "  -0.40%  "
$ws.Range("E38").Value = "  -0.57%  "
$ws.Range("E39").Value = "  -1.07%  "
$ws.Range("E40").Value = "  -1.55%  "
$ws.Range("E41").Value = "  -1.34%  "
$ws.Range("D42").Value = "'2.70"
$ws.Range("E42").Value = "  -2.82%  "
$ws.Range("D43").Value = "'2.009.11"
$ws.Range("E43").Value = "  -0.28%  "
$ws.Range("D45").Value = "'18.34"
$ws.Range("E45").Value = "  +4.62%  "
$ws.Range("D46").Value = "'10.07"
$ws.Range("E46").Value = "  -2.77%  "
$ws.Range("D47").Value = "'2.03"
$ws.Range("E47").Value = "  -8.23%  "
$ws.Range("E48").Value = "  -0.65%  "
$ws.Range("D49").Value = "'53.94"
$ws.Range("E49").Value = "  -0.30%  "
$ws.Range("D50").Value = "'2.85"
$ws.Range("E50").Value = "  +10.91%  "
$ws.Range("D51").Value = "'2.525.20"
$ws.Range("E51").Value = "  -0.07%  "
